$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data: date, task description, hours left blank
$ws.Range("A9").Value = 45623
$ws.Range("A9").NumberFormat = $ws.Range("A8").NumberFormat
$ws.Range("B9").Value = "Redoing backend for ui controllers so its easier to use."

# Update selection to match the new active cell (B9)
$ws.Range("B9").Select()
